$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Row 5
$ws.Range("A5").Value = 44201
$ws.Range("B5").Value = "'`$59.90"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'`$169.0"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'`$39.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'`$45.00"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'`$219.0"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'`$29.90"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = "'`$889.0"
$ws.Range("H5").Style = "Normal"
$ws.Range("I5").Value = "'`$719.0"
$ws.Range("I5").Style = "Normal"
$ws.Range("J5").Value = "'`$889.0"
$ws.Range("J5").Style = "Normal"
$ws.Range("K5").Value = "'`$369.0"
$ws.Range("K5").Style = "Normal"

# Row 6
$ws.Range("A6").Value = 44201
$ws.Range("B6").Value = "'`$59.90"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'`$169.0"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'`$39.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'`$45.00"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'`$219.0"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'`$29.90"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = "'`$889.0"
$ws.Range("H6").Style = "Normal"
$ws.Range("I6").Value = "'`$719.0"
$ws.Range("I6").Style = "Normal"
$ws.Range("J6").Value = "'`$889.0"
$ws.Range("J6").Style = "Normal"
$ws.Range("K6").Value = "'`$369.0"
$ws.Range("K6").Style = "Normal"

# Row 7
$ws.Range("A7").Value = 44203
$ws.Range("B7").Value = "'`$59.90"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'`$169.00"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'`$39.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'`$45.00"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'`$219.00"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'`$24.90`n`$29.90-17%"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("H7").Style = "Normal"
$ws.Range("I7").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("I7").Style = "Normal"
$ws.Range("J7").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("J7").Style = "Normal"
$ws.Range("K7").Value = "'`$289.00`n`$369.00-22%"
$ws.Range("K7").Style = "Normal"
$ws.Rows.Item(7).AutoFit()

# Row 8
$ws.Range("A8").Value = 44203
$ws.Range("B8").Value = "'`$59.90"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'`$169.00"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'`$39.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'`$45.00"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'`$219.00"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'`$24.90`n`$29.90-17%"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("H8").Style = "Normal"
$ws.Range("I8").Value = "'`$709.00`n`$719.00-1%"
$ws.Range("I8").Style = "Normal"
$ws.Range("J8").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("J8").Style = "Normal"
$ws.Range("K8").Value = "'`$289.00`n`$369.00-22%"
$ws.Range("K8").Style = "Normal"
$ws.Rows.Item(8).AutoFit()

# Row 9
$ws.Range("A9").Value = 44203
$ws.Range("B9").Value = "'`$59.90"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'`$169.00"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'`$39.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'`$45.00"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'`$219.00"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "'`$24.90`n`$29.90-17%"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("H9").Style = "Normal"
$ws.Range("I9").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("I9").Style = "Normal"
$ws.Range("J9").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("J9").Style = "Normal"
$ws.Range("K9").Value = "'`$289.00`n`$369.00-22%"
$ws.Range("K9").Style = "Normal"
$ws.Rows.Item(9).AutoFit()

# Row 10
$ws.Range("A10").Value = 44203
$ws.Range("B10").Value = "'`$59.90"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'`$169.00"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'`$39.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'`$45.00"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'`$219.00"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "'`$24.90`n`$29.90-17%"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("H10").Style = "Normal"
$ws.Range("I10").Value = "'`$709.00`n`$719.00-1%"
$ws.Range("I10").Style = "Normal"
$ws.Range("J10").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("J10").Style = "Normal"
$ws.Range("K10").Value = "'`$289.00`n`$369.00-22%"
$ws.Range("K10").Style = "Normal"
$ws.Rows.Item(10).AutoFit()

# Row 11
$ws.Range("A11").Value = 44203
$ws.Range("B11").Value = "'`$59.90"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'`$169.00"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'`$39.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'`$45.00"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'`$219.00"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "'`$24.90`n`$29.90-17%"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("H11").Style = "Normal"
$ws.Range("I11").Value = "'`$709.00`n`$719.00-1%"
$ws.Range("I11").Style = "Normal"
$ws.Range("J11").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("J11").Style = "Normal"
$ws.Range("K11").Value = "'`$289.00`n`$369.00-22%"
$ws.Range("K11").Style = "Normal"
$ws.Rows.Item(11).AutoFit()

# Row 12
$ws.Range("A12").Value = 44203
$ws.Range("A12").NumberFormat = "yyyy-mm-dd"
$ws.Range("B12").Value = "'`$59.90"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'`$169.00"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'`$39.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'`$45.00"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "'`$219.00"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = "'`$24.90`n`$29.90-17%"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("H12").Style = "Normal"
$ws.Range("I12").Value = "'`$709.00`n`$719.00-1%"
$ws.Range("I12").Style = "Normal"
$ws.Range("J12").Value = "'`$851.00`n`$889.00-4%"
$ws.Range("J12").Style = "Normal"
$ws.Range("K12").Value = "'`$289.00`n`$369.00-22%"
$ws.Range("K12").Style = "Normal"
$ws.Rows.Item(12).AutoFit()

